# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# for rows 2-51 per the latest scrape, applied via GitHub Actions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Assign a value while forcing Excel to keep it as plain text,
    # even if the text looks like a number (e.g. '231.07') or has
    # multiple dots (e.g. '34.657.21'). The leading apostrophe forces
    # text interpretation; resetting the style back to Normal removes
    # the quote-prefix formatting Excel applies, leaving the cell's
    # style untouched while keeping the exact string value.
    $cell.Value = "'" + $text
    $cell.Style = 'Normal'
}

Set-TextValue $ws.Range('D2') '34.657.21'
Set-TextValue $ws.Range('E2') '  -2.22%  '
Set-TextValue $ws.Range('E3') '  -2.35%  '
Set-TextValue $ws.Range('E4') '  +0.10%  '
Set-TextValue $ws.Range('D5') '231.07'
Set-TextValue $ws.Range('E5') '  +0.59%  '
Set-TextValue $ws.Range('E6') '  -1.06%  '
Set-TextValue $ws.Range('E7') '  +0.09%  '
Set-TextValue $ws.Range('D8') '39.06'
Set-TextValue $ws.Range('E8') '  -7.20%  '
Set-TextValue $ws.Range('D9') '0.317'
Set-TextValue $ws.Range('E9') '  +3.57%  '
Set-TextValue $ws.Range('E10') '  -1.96%  '
Set-TextValue $ws.Range('D11') '0.0990'
Set-TextValue $ws.Range('E11') '  -1.88%  '
Set-TextValue $ws.Range('D12') '2.066.21'
Set-TextValue $ws.Range('E12') '  -2.42%  '
Set-TextValue $ws.Range('D13') '1.808.72'
Set-TextValue $ws.Range('E13') '  -2.23%  '
Set-TextValue $ws.Range('D14') '0.660'
Set-TextValue $ws.Range('E14') '  -1.62%  '
Set-TextValue $ws.Range('D15') '10.87'
Set-TextValue $ws.Range('E15') '  -4.75%  '
Set-TextValue $ws.Range('D16') '4.54'
Set-TextValue $ws.Range('E16') '  -2.98%  '
Set-TextValue $ws.Range('D17') '34.666.72'
Set-TextValue $ws.Range('E17') '  -2.22%  '
Set-TextValue $ws.Range('D18') '69.45'
Set-TextValue $ws.Range('E18') '  -0.80%  '
Set-TextValue $ws.Range('D19') '0.0₃0781'
Set-TextValue $ws.Range('D20') '238.91'
Set-TextValue $ws.Range('E20') '  -3.40%  '
Set-TextValue $ws.Range('D21') '11.75'
Set-TextValue $ws.Range('E21') '  -2.69%  '
Set-TextValue $ws.Range('D22') '4.65'
Set-TextValue $ws.Range('E22') '  +0.30%  '
Set-TextValue $ws.Range('E23') '  +0.10%  '
Set-TextValue $ws.Range('D24') '2.23'
Set-TextValue $ws.Range('E24') '  +1.85%  '
Set-TextValue $ws.Range('D25') '172.53'
Set-TextValue $ws.Range('E25') '  +1.96%  '
Set-TextValue $ws.Range('D26') '7.68'
Set-TextValue $ws.Range('E26') '  -3.03%  '
Set-TextValue $ws.Range('D27') '17.12'
Set-TextValue $ws.Range('E28') '  -2.49%  '
Set-TextValue $ws.Range('D29') '1.52'
Set-TextValue $ws.Range('E29') '  +9.13%  '
Set-TextValue $ws.Range('E30') '  +0.16%  '
Set-TextValue $ws.Range('D31') '3.98'
Set-TextValue $ws.Range('E31') '  +1.36%  '
Set-TextValue $ws.Range('D32') '0.0543'
Set-TextValue $ws.Range('E32') '  -0.13%  '
Set-TextValue $ws.Range('D33') '3.93'
Set-TextValue $ws.Range('E33') '  -3.03%  '
Set-TextValue $ws.Range('D34') '1.26'
Set-TextValue $ws.Range('E34') '  +16.14%  '
Set-TextValue $ws.Range('D35') '1.76'
Set-TextValue $ws.Range('E35') '  -5.70%  '
Set-TextValue $ws.Range('D36') '0.690'
Set-TextValue $ws.Range('E36') '  +0.76%  '
Set-TextValue $ws.Range('D37') '90.81'
Set-TextValue $ws.Range('E37') '  -7.00%  '
Set-TextValue $ws.Range('E38') '  +4.98%  '
Set-TextValue $ws.Range('D39') '1.308.60'
Set-TextValue $ws.Range('E39') '  -3.41%  '
Set-TextValue $ws.Range('E40') '  -1.95%  '
Set-TextValue $ws.Range('E41') '  -0.75%  '
Set-TextValue $ws.Range('D42') '0.956'
Set-TextValue $ws.Range('E42') '  -4.36%  '
Set-TextValue $ws.Range('D43') '14.17'
Set-TextValue $ws.Range('E43') '  -3.23%  '
Set-TextValue $ws.Range('D44') '2.21'
Set-TextValue $ws.Range('E44') '  -10.04%  '
Set-TextValue $ws.Range('E45') '  -5.60%  '
Set-TextValue $ws.Range('D46') '6.13'
Set-TextValue $ws.Range('E46') '  -1.08%  '
Set-TextValue $ws.Range('D47') '0.0509'
Set-TextValue $ws.Range('E47') '  -1.93%  '
Set-TextValue $ws.Range('D48') '1.994.30'
Set-TextValue $ws.Range('E48') '  -1.09%  '
Set-TextValue $ws.Range('E49') '  +0.09%  '
Set-TextValue $ws.Range('E50') '  +7.87%  '
Set-TextValue $ws.Range('D51') '98.53'
Set-TextValue $ws.Range('E51') '  -4.89%  '
